# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-string suffixes to "_FV2404" / "_FV2410"
# - Freeze the header row (row 1) in the sheet view
# - Wrap the data range A1:U55 in an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row shared strings -----------------------------------
# Columns A-J carried the "_old" suffix -> becomes "_FV2404"
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
# Columns L-U carried the "_new" suffix -> becomes "_FV2410" (column K = "diff" is untouched)
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $fv2404Headers[$i]
}

for ($i = 0; $i -lt $fv2410Headers.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $fv2410Headers[$i]
}

# --- 2. Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Turn the used range into an Excel table -----------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), $null, 1)
$tbl.Name = "Table1"
